$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21
$ws.Range("A7").Value = -21.255
$ws.Range("B7").Value = 6.536
$ws.Range("A8").Value = -21.623
$ws.Range("B11").Value = 6.218
$ws.Range("B12").Value = 5.414999999999999
$ws.Range("D12").Value = -7.105
$ws.Range("D13").Value = -8.103999999999999
$ws.Range("D14").Value = -7.661
$ws.Range("B15").Value = 5.359
$ws.Range("A16").Value = -21.479
$ws.Range("D16").Value = -8.518000000000001
$ws.Range("D19").Value = -7.796000000000001
$ws.Range("A20").Value = -22.103
$ws.Range("B20").Value = 5.928
$ws.Range("D20").Value = -7.965999999999999
$ws.Range("A21").Value = -20.952
$ws.Range("B21").Value = 7.692000000000002
$ws.Range("B22").Value = 6.391000000000001
$ws.Range("D22").Value = -8.151
$ws.Range("B23").Value = 7.359999999999999
$ws.Range("A28").Value = -21.778
$ws.Range("A29").Value = -21.675
$ws.Range("B29").Value = 5.816999999999999
$ws.Range("A30").Value = -21.744
$ws.Range("A32").Value = -21.7
$ws.Range("B34").Value = 8.059000000000001
$ws.Range("D36").Value = -7.834000000000001
$ws.Range("A40").Value = -20.605
$ws.Range("B42").Value = 7.047000000000001
$ws.Range("B43").Value = 5.744
$ws.Range("D43").Value = -8.099
$ws.Range("B44").Value = 5.34
$ws.Range("B45").Value = 5.403999999999999
$ws.Range("A46").Value = -20.849
$ws.Range("B46").Value = 6.865
$ws.Range("D46").Value = -8.138
$ws.Range("B50").Value = 5.637
$ws.Range("D50").Value = -7.897
$ws.Range("A51").Value = -20.952
$ws.Range("B51").Value = 7.930000000000001
$ws.Range("A52").Value = -21.588
$ws.Range("A57").Value = -21.614
$ws.Range("B57").Value = 6.078
$ws.Range("A59").Value = -22.091
$ws.Range("A62").Value = -22.028
$ws.Range("B65").Value = 5.306
$ws.Range("A66").Value = -21.47
$ws.Range("B66").Value = 5.709000000000001
$ws.Range("B67").Value = 5.778
$ws.Range("A73").Value = -20.618
$ws.Range("A74").Value = -21.043
$ws.Range("D76").Value = -7.993
$ws.Range("A77").Value = -21.298
$ws.Range("B79").Value = 5.680999999999999
$ws.Range("B84").Value = 5.781000000000001
$ws.Range("B87").Value = 4.447000000000001
$ws.Range("A92").Value = -21.511
$ws.Range("B92").Value = 5.527
$ws.Range("D95").Value = -7.952000000000001
$ws.Range("B97").Value = 5.161
$ws.Range("D97").Value = -8.465
$ws.Range("D99").Value = -7.772
$ws.Range("A100").Value = -21.481
